$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E4"  = 16.45259999999999
    "E6"  = 16.53229999999999
    "E7"  = 15.6725
    "E16" = 15.9293
    "E20" = 15.91219999999999
    "E28" = 16.60469999999999
    "E29" = 17.23260000000002
    "E32" = 16.95759999999999
    "E40" = 17.0403
    "E46" = 17.11569999999999
    "E51" = 17.2709
    "E52" = 17.0797
    "E57" = 16.6178
    "E59" = 15.79510000000001
    "E62" = 16.6909
    "E66" = 17.11600000000001
    "E73" = 17.04190000000001
    "E74" = 16.66079999999998
    "E92" = 18.49360000000003
    "E100" = 16.4122
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
